# Commit: "Using internal column widths in pptx writer tables (#9392)"
#
# The pptx table writer used to always split the available width evenly
# across all columns. Now it prefers the incoming (explicit) column
# widths when present, and only falls back to even distribution when
# they are missing. Because of the different rounding this produces,
# a table whose two columns were previously written as 2501900 EMU
# (197 pt) each now gets 2514600 EMU (198 pt) each - the same width
# already used by the other (untouched) table on the slide.
#
# EMU/point: 12700 EMU = 1 pt -> 2501900 EMU = 197 pt, 2514600 EMU = 198 pt.

$oldWidthPt = 197
$newWidthPt = 198
$tolerance = 0.01

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $col = $tbl.Columns.Item($c)
                if ([Math]::Abs($col.Width - $oldWidthPt) -lt $tolerance) {
                    $col.Width = $newWidthPt
                }
            }
        }
    }
}
